{"js": "// Update the three-digit x one-digit multiplication prompts to the\n// newly generated set of problems. Each cell's text is a unique\n// \"NNN\u00d7N=\" string, so we can safely match-and-replace each one.\nconst replacements = [\n  [\"532\u00d77=\", \"237\u00d76=\"],\n  [\"773\u00d72=\", \"650\u00d79=\"],\n  [\"720\u00d72=\", \"710\u00d74=\"],\n  [\"444\u00d78=\", \"877\u00d72=\"],\n  [\"848\u00d74=\", \"794\u00d74=\"],\n  [\"106\u00d72=\", \"304\u00d79=\"],\n  [\"848\u00d73=\", \"651\u00d78=\"],\n  [\"556\u00d72=\", \"240\u00d77=\"],\n  [\"615\u00d75=\", \"458\u00d74=\"],\n  [\"856\u00d78=\", \"758\u00d75=\"],\n  [\"384\u00d74=\", \"703\u00d75=\"],\n  [\"293\u00d77=\", \"525\u00d78=\"],\n  [\"771\u00d75=\", \"558\u00d77=\"],\n  [\"611\u00d79=\", \"359\u00d73=\"],\n  [\"949\u00d74=\", \"906\u00d76=\"],\n  [\"226\u00d78=\", \"561\u00d79=\"],\n  [\"431\u00d75=\", \"926\u00d74=\"],\n  [\"847\u00d78=\", \"773\u00d73=\"],\n  [\"977\u00d78=\", \"614\u00d72=\"],\n  [\"571\u00d74=\", \"833\u00d73=\"],\n  [\"898\u00d75=\", \"922\u00d76=\"],\n  [\"923\u00d79=\", \"657\u00d73=\"],\n  [\"391\u00d77=\", \"978\u00d77=\"],\n  [\"951\u00d79=\", \"116\u00d78=\"],\n  [\"426\u00d74=\", \"792\u00d74=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the three-digit x one-digit multiplication prompts to the\n# newly generated set of problems. Each cell's text is a unique\n# \"NNN\u00d7N=\" string, so Find/Replace on the exact string is unambiguous.\n$pairs = @(\n  @(\"532\u00d77=\", \"237\u00d76=\"),\n  @(\"773\u00d72=\", \"650\u00d79=\"),\n  @(\"720\u00d72=\", \"710\u00d74=\"),\n  @(\"444\u00d78=\", \"877\u00d72=\"),\n  @(\"848\u00d74=\", \"794\u00d74=\"),\n  @(\"106\u00d72=\", \"304\u00d79=\"),\n  @(\"848\u00d73=\", \"651\u00d78=\"),\n  @(\"556\u00d72=\", \"240\u00d77=\"),\n  @(\"615\u00d75=\", \"458\u00d74=\"),\n  @(\"856\u00d78=\", \"758\u00d75=\"),\n  @(\"384\u00d74=\", \"703\u00d75=\"),\n  @(\"293\u00d77=\", \"525\u00d78=\"),\n  @(\"771\u00d75=\", \"558\u00d77=\"),\n  @(\"611\u00d79=\", \"359\u00d73=\"),\n  @(\"949\u00d74=\", \"906\u00d76=\"),\n  @(\"226\u00d78=\", \"561\u00d79=\"),\n  @(\"431\u00d75=\", \"926\u00d74=\"),\n  @(\"847\u00d78=\", \"773\u00d73=\"),\n  @(\"977\u00d78=\", \"614\u00d72=\"),\n  @(\"571\u00d74=\", \"833\u00d73=\"),\n  @(\"898\u00d75=\", \"922\u00d76=\"),\n  @(\"923\u00d79=\", \"657\u00d73=\"),\n  @(\"391\u00d77=\", \"978\u00d77=\"),\n  @(\"951\u00d79=\", \"116\u00d78=\"),\n  @(\"426\u00d74=\", \"792\u00d74=\")\n)\n\n$d = $word.ActiveDocument\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Text = $oldText\n  $find.Replacement.ClearFormatting()\n  $find.Replacement.Text = $newText\n  $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
